$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newAuthors = "[Ke%Hu%NULL%1,                  Wei-jie%Guan%NULL%1,                  Ying%Bi%NULL%1,                  Wei%Zhang%NULL%0,                  Lanjuan%Li%NULL%0,                  Boli%Zhang%NULL%1,                  Qingquan%Liu%NULL%1,                  Yuanlin%Song%NULL%1,                  Xingwang%Li%NULL%0,                  Zhongping%Duan%NULL%1,                  Qingshan%Zheng%NULL%1,                  Zifeng%Yang%NULL%1,                  Jingyi%Liang%NULL%1,                  Mingfeng%Han%NULL%0,                  Lianguo%Ruan%NULL%1,                  Chaomin%Wu%NULL%1,                  Yunting%Zhang%NULL%1,                  Zhen-hua%Jia%NULL%1,                  Nan-shan%Zhong%NULL%0]"

$ws.Range("E2").Value = $newAuthors
$ws.Range("I2").Value = "_elsevier_CROSSREF"
